$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New flag string used by several rows below
$flagText = "Flagged: Low risk but took less than 1 month."

# Row 2: date shifts to 2024-06-01, stays "Created as Case Accepted", count 0
$ws.Range("A2").Value = 45444
$ws.Range("C2").Value = 0

# Row 3: date shifts to 2024-06-01, becomes flagged, count 30
$ws.Range("A3").Value = 45444
$ws.Range("B3").Value = $flagText
$ws.Range("C3").Value = 30

# Row 4: date shifts to 2024-07-01, stays "Created as Case Accepted", count 0
$ws.Range("A4").Value = 45474
$ws.Range("C4").Value = 0

# Row 5: date shifts to 2024-07-01, becomes flagged, count 0
$ws.Range("A5").Value = 45474
$ws.Range("B5").Value = $flagText
$ws.Range("C5").Value = 0

# Row 6: date shifts to 2024-08-01, stays "Created as Case Accepted", count 0
$ws.Range("A6").Value = 45505
$ws.Range("C6").Value = 0

# Row 7: date shifts to 2024-08-01, becomes flagged, count 0
$ws.Range("A7").Value = 45505
$ws.Range("B7").Value = $flagText
$ws.Range("C7").Value = 0

# Row 8: date shifts to 2024-09-01, stays "Created as Case Accepted", count 63
$ws.Range("A8").Value = 45536
$ws.Range("C8").Value = 63

# Row 9: date shifts to 2024-09-01, becomes flagged, count 0
$ws.Range("A9").Value = 45536
$ws.Range("B9").Value = $flagText
$ws.Range("C9").Value = 0

# Drop the former rows 10-13 entirely
$ws.Rows("10:13").Delete()
